$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.144.83"
$ws.Range("E2").Value = "  -0.03%  "

$ws.Range("D3").Value = "3.128.43"
$ws.Range("E3").Value = "  +0.16%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.11%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.10"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.47%  "

$ws.Range("E7").Value = "  +0.10%  "

$ws.Range("D8").Value = "3.125.53"
$ws.Range("E8").Value = "  +0.23%  "

$ws.Range("E9").Value = "  -1.13%  "

$ws.Range("E10").Value = "  -0.60%  "

$ws.Range("E11").Value = "  -1.96%  "

$ws.Range("E12").Value = "  -0.73%  "

$ws.Range("E13").Value = "  -2.58%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.44"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.95%  "

$ws.Range("E15").Value = "  -0.32%  "

$ws.Range("D16").Value = "3.651.08"
$ws.Range("E16").Value = "  +0.26%  "

$ws.Range("D17").Value = "67.083.67"
$ws.Range("E17").Value = "  -0.10%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.05"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.09%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.09"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.16%  "

$ws.Range("D20").Value = "3.129.29"
$ws.Range("E20").Value = "  +0.21%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "490.86"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.12%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.81"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.50%  "

$ws.Range("E23").Value = "  -1.69%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.98"
$ws.Range("D24").Style = "Normal"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.83"
$ws.Range("D25").Style = "Normal"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.28"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.32%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.30"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.32%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.08%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.05"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.39%  "

$ws.Range("E30").Value = "  -1.93%  "

$ws.Range("E31").Value = "  -2.48%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.21"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.40%  "

$ws.Range("E33").Value = "  -1.11%  "

$ws.Range("E34").Value = "  -0.46%  "

$ws.Range("E35").Value = "  +0.04%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "48.38"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.23%  "

$ws.Range("E37").Value = "  -3.75%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.947"
$ws.Range("D38").Style = "Normal"

$ws.Range("E39").Value = "  -1.32%  "

$ws.Range("E40").Value = "  -0.03%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.02"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.92%  "

$ws.Range("E42").Value = "  -0.52%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.70"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.28%  "

$ws.Range("D45").Value = "2.803.88"
$ws.Range("E45").Value = "  -0.76%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "376.24"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.47%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0348"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.54%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "135.13"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.29%  "

$ws.Range("E49").Value = "  +0.00%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.01"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.52%  "

$ws.Range("E51").Value = "  +2.18%  "

